$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "54.697.43"
$ws.Range("E2").Value = "  +5.76%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.189.49"
$ws.Range("E3").Value = "  +2.95%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "402.04"
$ws.Range("E5").Value = "  +3.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "108.73"
$ws.Range("E6").Value = "  +5.15%  "

$ws.Range("E7").Value = "  +1.34%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.05%  "

$ws.Range("E9").Value = "  +5.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.07"
$ws.Range("E10").Value = "  +5.51%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.679.14"
$ws.Range("E13").Value = "  +2.67%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.05"
$ws.Range("E14").Value = "  +1.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.05"
$ws.Range("E15").Value = "  +3.16%  "

$ws.Range("E16").Value = "  +8.87%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.189.68"
$ws.Range("E17").Value = "  +3.14%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "10.56"
$ws.Range("E18").Value = "  -0.93%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "54.530.67"
$ws.Range("E19").Value = "  +5.15%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.33"
$ws.Range("E20").Value = "  +3.61%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.91"
$ws.Range("E21").Value = "  +3.45%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0000100"
$ws.Range("E22").Value = "  +3.42%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.74"
$ws.Range("E23").Value = "  +3.88%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "275.60"
$ws.Range("E24").Value = "  +2.51%  "

$ws.Range("E25").Value = "  +4.87%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.07"
$ws.Range("E26").Value = "  -0.89%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.83"
$ws.Range("E27").Value = "  +2.60%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.46"
$ws.Range("E28").Value = "  +3.65%  "

$ws.Range("E29").Value = "  -0.70%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  -0.11%  "

$ws.Range("E31").Value = "  +4.22%  "

$ws.Range("E32").Value = "  +6.67%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0506"
$ws.Range("E33").Value = "  +12.58%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "37.00"
$ws.Range("E34").Value = "  +3.39%  "

$ws.Range("E35").Value = "  +0.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.80"
$ws.Range("E36").Value = "  +0.64%  "

$ws.Range("E37").Value = "  +7.88%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("E38").Value = "  -0.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.85"
$ws.Range("E39").Value = "  +9.83%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.18"
$ws.Range("E40").Value = "  +13.00%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.95"
$ws.Range("E41").Value = "  +3.37%  "

$ws.Range("E42").Value = "  -0.06%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.30"
$ws.Range("E43").Value = "  +1.90%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "130.94"
$ws.Range("E44").Value = "  +2.60%  "

$ws.Range("E45").Value = "  +1.58%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.33"
$ws.Range("E46").Value = "  +0.98%  "

$ws.Range("E47").Value = "  -1.23%  "

$ws.Range("E48").Value = "  -1.03%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.091.88"
$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0350"
$ws.Range("E50").Value = "  +9.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0511"
$ws.Range("E51").Value = "  +11.19%  "

# Row 11/12 coin swap (TRON <-> Dogecoin)
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0887"
$ws.Range("E11").Value = "  +2.94%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.140"
$ws.Range("E12").Value = "  +1.59%  "